# Update "想去人数" (want-to-go count) figures on the "展览" and "全部类型" sheets
# to reflect freshly generated data (gh-pages output at 456a3b4).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 348   # was 341
$ws1.Range("F4").Value = 4680  # was 4660
$ws1.Range("F5").Value = 44    # was 45
$ws1.Range("F6").Value = 473   # was 471

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 348   # was 341
$ws4.Range("F4").Value = 4680  # was 4660
$ws4.Range("F7").Value = 44    # was 45
$ws4.Range("F8").Value = 473   # was 471

$wb.Save()
